# Add "2022-Q1" sheet after "2021-Q4" and before "总计"; populate its fund-holdings
# table; and update the "总计" summary sheet with a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" worksheet, positioned right after "2021-Q4".
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "2022-Q1"

# Borrow header/column-A formatting (bold, centered, bordered) from "2021-Q4" so the
# new sheet matches the look of the other quarterly sheets.
$afterSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$afterSheet.Range("A2").Copy()
$newSheet.Range("A2:A33").PasteSpecial(-4122)

# Header row
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

$fundRows = @(
    @("011333","鹏华品质优选混合A","44.36","88.31","6.64","2.9455",5),
    @("011570","鹏华鑫远价值一年持有期混合型证券投资基金A","12.40","85.38","8.44","1.0466",4),
    @("012057","鹏华品质成长混合型证券投资基金A","13.38","81.87","6.74","0.9018",6),
    @("010591","富国中国中小盘混合(QDII)美元","35.75","86.53","2.41","0.8616",9),
    @("100061","富国中国中小盘混合(QDII)人民币","35.75","86.53","2.41","0.8616",9),
    @("009984","鹏华启航两年封闭运作混合","14.73","89.97","4.38","0.6452",3),
    @("009234","鹏华优质企业混合","4.14","91.73","7.26","0.3006",6),
    @("009223","宝盈现代服务业混合A","3.63","91.38","6.60","0.2396",7),
    @("013859","宝盈品质甄选混合A","3.26","92.14","6.28","0.2047",7),
    @("011334","鹏华品质优选混合C","2.84","88.31","6.64","0.1886",5),
    @("007548","易方达ESG责任投资股票","2.92","86.09","4.13","0.1206",10),
    @("501021","华宝兴业标普香港上市中国中小盘指数(QDII-LOF)A","4.96","94.77","1.87","0.0928",7),
    @("011571","鹏华鑫远价值一年持有期混合型证券投资基金C","0.78","85.38","8.44","0.0658",4),
    @("012058","鹏华品质成长混合型证券投资基金C","0.47","81.87","6.74","0.0317",6),
    @("009224","宝盈现代服务业混合C","0.46","91.38","6.60","0.0304",7),
    @("005142","中融沪港深大消费主题灵活配置混合A","0.49","88.98","6.08","0.0298",1),
    @("006675","宝盈品牌消费股票A","0.32","90.31","6.42","0.0205",8),
    @("005143","中融沪港深大消费主题灵活配置混合C","0.33","88.98","6.08","0.0201",1),
    @("007751","景顺长城中证沪港深红利成长低波动指数A","0.83","91.29","2.14","0.0178",9),
    @("004532","民生加银中证港股通高股息精选指数A","0.26","94.88","4.82","0.0125",4),
    @("006676","宝盈品牌消费股票C","0.15","90.31","6.42","0.0096",8),
    @("011647","博时港股通红利精选混合A","0.13","92.10","5.69","0.0074",5),
    @("501303","广发港股通恒生综合中型股指数(LOF)A","0.34","92.39","1.90","0.0065",4),
    @("004533","民生加银中证港股通高股息精选指数C","0.10","94.88","4.82","0.0048",4),
    @("006127","华宝兴业标普香港上市中国中小盘指数(QDII-LOF)C","0.23","94.77","1.87","0.0043",7),
    @("001942","前海开源沪港深汇鑫灵活配置混合A","0.10","90.39","3.25","0.0032",9),
    @("013860","宝盈品质甄选混合C","0.05","92.14","6.28","0.0031",7),
    @("001943","前海开源沪港深汇鑫灵活配置混合C","0.08","90.39","3.25","0.0026",9),
    @("004996","广发港股通恒生综合中型股指数(LOF)C","0.11","92.39","1.90","0.0021",4),
    @("160922","大成恒生综合中小型股指数(QDII-LOF)A","0.10","92.44","1.44","0.0014",5),
    @("007760","景顺长城中证沪港深红利成长低波动指数C","0.06","91.29","2.14","0.0013",9),
    @("011648","博时港股通红利精选混合C","0.02","92.10","5.69","0.0011",5)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $newSheet.Cells.Item($r, 1).Value = $i

    $newSheet.Cells.Item($r, 2).NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $row[0]

    $newSheet.Cells.Item($r, 3).NumberFormat = "@"
    $newSheet.Cells.Item($r, 3).Value = $row[1]

    $newSheet.Cells.Item($r, 4).NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $row[2]

    $newSheet.Cells.Item($r, 5).NumberFormat = "@"
    $newSheet.Cells.Item($r, 5).Value = $row[3]

    $newSheet.Cells.Item($r, 6).NumberFormat = "@"
    $newSheet.Cells.Item($r, 6).Value = $row[4]

    $newSheet.Cells.Item($r, 7).NumberFormat = "@"
    $newSheet.Cells.Item($r, 7).Value = $row[5]

    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row 2 for "2022-Q1" and
#    renumber the existing rows' index column (A) by +1.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The inserted row inherits row 1's (header) bold/centered format on B:D; strip it
# so it matches the plain look of the other data rows, then restore column A's
# usual bold/centered/bordered index style by copying it from row 3.
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 32
$totalSheet.Cells.Item(2,4).Value = 8.69

for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
